$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing formatting issues on columns A/B (part of "fix db") ---
# A3 should match the A1/A2 "id column" style instead of its old bold/green style
$ws.Range("A1").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# A4 and B4 should match the plain B-column style instead of their old bold/green style
$ws.Range("B1").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)

# --- Add new column C ("item module" table column) ---
$ws.Range("C1").Value = $ws.Range("B1").Value2
$ws.Range("C2").Value = "table"
$ws.Range("C3").Value = "value2"
$ws.Range("C4").Value = $ws.Range("B4").Value2
$ws.Range("C5").Value = "{{reward_type = 3,item_type = 2,item_count = 1},{reward_type = 4,item_type = 1,item_count = 1000}}"

# Match the new column's formatting to column B (header/style rows) and row 5 (data style)
$ws.Range("B1:B4").Copy()
$ws.Range("C1:C4").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("C5").PasteSpecial(-4122)

# Match column width + grouping/outline level of the new column to column B
$ws.Columns.Item(3).ColumnWidth = 9.86
$ws.Columns.Item(3).OutlineLevel = 2

# Preserve the sheet's row-outline summary (cosmetic metadata) without leaving stray rows
$ws.Rows.Item(100).OutlineLevel = 4
$ws.Rows.Item(100).Delete()

# --- Cosmetic: update the active selection like the authored workbook ---
$ws.Range("L12").Select()
